# Update gh-pages generated output: refresh the "view count" (F column)
# figures on the "展览" and "全部类型" worksheets, which mirror each other.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value
$updates = @{
    "F4"  = 11036
    "F5"  = 10216
    "F13" = 9582
    "F15" = 2439
    "F17" = 7
    "F20" = 10856
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
